# Update the "想去人数" (want-to-go count) values in column F
# for the "展览" sheet and the "全部类型" sheet to reflect the
# latest scraped counts (gh-pages output regeneration).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 1370
$ws1.Range("F3").Value  = 1621
$ws1.Range("F7").Value  = 680
$ws1.Range("F11").Value = 2478
$ws1.Range("F13").Value = 1510
$ws1.Range("F15").Value = 249
$ws1.Range("F18").Value = 85
$ws1.Range("F19").Value = 313
$ws1.Range("F24").Value = 5103
$ws1.Range("F26").Value = 566
$ws1.Range("F27").Value = 86
$ws1.Range("F34").Value = 745
$ws1.Range("F39").Value = 1080
$ws1.Range("F42").Value = 176

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 1370
$ws4.Range("F5").Value  = 1621
$ws4.Range("F11").Value = 680
$ws4.Range("F17").Value = 2478
$ws4.Range("F19").Value = 1510
$ws4.Range("F21").Value = 249
$ws4.Range("F25").Value = 85
$ws4.Range("F26").Value = 313
$ws4.Range("F29").Value = 5103
$ws4.Range("F31").Value = 566
$ws4.Range("F32").Value = 86
$ws4.Range("F39").Value = 745
$ws4.Range("F42").Value = 1080
$ws4.Range("F44").Value = 176
